$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Section "Medidas ponderadas" (row 29) -> "Tiempo de respuesta" weighted calc
# ---------------------------------------------------------------------------
$ws.Range("A29").Value = "Medidas ponderadas"
$ws.Range("A30").Value = "Tiempo de respuesta"

# Header row (row 31) - copy style from the existing bordered header row (row 2)
$ws.Range("B2:I2").Copy($ws.Range("B31:I31"))
$ws.Range("B31").Value = "build-mplayer"
$ws.Range("C31").Value = "build-php"
$ws.Range("D31").Value = "compress-gzip"
$ws.Range("E31").Value = "dcraw"
$ws.Range("F31").Value = "encode-flac"
$ws.Range("G31").Value = "gnupg"
$ws.Range("H31").Value = "mafft"
$ws.Range("I31").Value = "mrbayes"

# Pesos row (row 32)
$ws.Range("A32").Value = "Pesos"
$ws.Range("B32").Value = 0.25
$ws.Range("C32").Value = 0.25
$ws.Range("D32").Value = 0.1
$ws.Range("E32").Value = 0.1
$ws.Range("F32").Value = 0.1
$ws.Range("G32").Value = 0.1
$ws.Range("H32").Value = 0.05
$ws.Range("I32").Value = 0.25

# Computadora (inv) block (rows 34-38)
$ws.Range("A34").Value = "Computadora (inv)"

$ws.Range("A35").Value = "A"
$ws.Range("B35").Formula = "=1/B3"
$ws.Range("C35").Formula = "=1/C3"
$ws.Range("D35").Formula = "=1/D3"
$ws.Range("E35").Formula = "=1/E3"
$ws.Range("F35").Formula = "=1/F3"
$ws.Range("G35").Formula = "=1/G3"
$ws.Range("H35").Formula = "=1/H3"
$ws.Range("I35").Formula = "=1/I3"

$ws.Range("A36").Value = "B"
$ws.Range("B36").Formula = "=1/B4"
$ws.Range("C36").Formula = "=1/C4"
$ws.Range("D36").Formula = "=1/D4"
$ws.Range("E36").Formula = "=1/E4"
$ws.Range("F36").Formula = "=1/F4"
$ws.Range("G36").Formula = "=1/G4"
$ws.Range("H36").Formula = "=1/H4"
$ws.Range("I36").Formula = "=1/I4"

$ws.Range("A37").Value = "C"
$ws.Range("B37").Formula = "=1/B5"
$ws.Range("C37").Formula = "=1/C5"
$ws.Range("D37").Formula = "=1/D5"
$ws.Range("E37").Formula = "=1/E5"
$ws.Range("F37").Formula = "=1/F5"
$ws.Range("G37").Formula = "=1/G5"
$ws.Range("H37").Formula = "=1/H5"
$ws.Range("I37").Formula = "=1/I5"

$ws.Range("A38").Value = "D"
$ws.Range("B38").Formula = "=1/B6"
$ws.Range("C38").Formula = "=1/C6"
$ws.Range("D38").Formula = "=1/D6"
$ws.Range("E38").Formula = "=1/E6"
$ws.Range("F38").Formula = "=1/F6"
$ws.Range("G38").Formula = "=1/G6"
$ws.Range("H38").Formula = "=1/H6"
$ws.Range("I38").Formula = "=1/I6"

# "Computadora" / "Media armónica" block (rows 43-47)
$ws.Range("A43").Value = "Computadora"
$ws.Range("B43").Value = "Media armónica"

# Row labels A44:A47 reuse the bordered style from A3:A6
$ws.Range("A3").Copy($ws.Range("A44"))
$ws.Range("A4").Copy($ws.Range("A45"))
$ws.Range("A5").Copy($ws.Range("A46"))
$ws.Range("A6").Copy($ws.Range("A47"))
$ws.Range("A44").Value = "A"
$ws.Range("A45").Value = "B"
$ws.Range("A46").Value = "C"
$ws.Range("A47").Value = "D"

$ws.Range("B44").Formula = "=SUM(B32:I32)/SUMPRODUCT(B32:I32, B35:I35)"
$ws.Range("B45").Formula = "=SUM(B32:I32)/SUMPRODUCT(B32:I32, B36:I36)"
$ws.Range("B46").Formula = "=SUM(B32:I32)/SUMPRODUCT(B32:I32, B37:I37)"
$ws.Range("B47").Formula = "=SUM(B32:I32)/SUMPRODUCT(B32:I32, B38:I38)"

# ---------------------------------------------------------------------------
# Section "Desempeño" (row 49) -> weighted Pc/Tareas calc
# ---------------------------------------------------------------------------
$ws.Range("A49").Value = "Desempeño"

$ws.Range("B23:F23").Copy($ws.Range("B50:F50"))
$ws.Range("B50").Value = "redis(LPOP)"
$ws.Range("C50").Value = "redis(SADD)"
$ws.Range("D50").Value = "redis(LPUSH)"
$ws.Range("E50").Value = "redis(GET)"
$ws.Range("F50").Value = "redis(SET)"

$ws.Range("A51").Value = "Pesos"
$ws.Range("B51").Value = 0.25
$ws.Range("C51").Value = 0.15
$ws.Range("D51").Value = 0.15
$ws.Range("E51").Value = 0.15
$ws.Range("F51").Value = 0.3

$ws.Range("A53").Value = "Computadora"
$ws.Range("B53").Value = "Media aritmética"

$ws.Range("A9").Copy($ws.Range("A54"))
$ws.Range("A10").Copy($ws.Range("A55"))
$ws.Range("A11").Copy($ws.Range("A56"))
$ws.Range("A12").Copy($ws.Range("A57"))
$ws.Range("A54").Value = "A"
$ws.Range("A55").Value = "B"
$ws.Range("A56").Value = "C"
$ws.Range("A57").Value = "D"

$ws.Range("B54").Formula = "=SUMPRODUCT(B9:F9, B51:F51)"
$ws.Range("B55").Formula = "=SUMPRODUCT(B10:F10, B51:F51)"
$ws.Range("B56").Formula = "=SUMPRODUCT(B11:F11, B51:F51)"
$ws.Range("B57").Formula = "=SUMPRODUCT(B12:F12, B51:F51)"

# ---------------------------------------------------------------------------
# Cosmetic tweaks: empty spacer rows, column width + selection
# ---------------------------------------------------------------------------
$ws.Rows.Item(33).RowHeight = 12.8
$ws.Rows.Item(39).RowHeight = 12.8
$ws.Rows.Item(40).RowHeight = 12.8
$ws.Rows.Item(41).RowHeight = 12.8
$ws.Rows.Item(42).RowHeight = 12.8
$ws.Rows.Item(48).RowHeight = 12.8
$ws.Rows.Item(52).RowHeight = 12.8

$ws.Columns.Item(1).ColumnWidth = 15.3
$ws.Range("B44:B47").Select() | Out-Null
